$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 1666.6666
$ws.Cells.Item(69, 10).Value = 1569.2307
$ws.Cells.Item(69, 12).Value = 4707.6921
$ws.Cells.Item(69, 14).Value = -6455.6921
$ws.Cells.Item(72, 8).Value = 1666.6666
$ws.Cells.Item(72, 10).Value = 1569.2307
$ws.Cells.Item(72, 12).Value = 14123.0763
$ws.Cells.Item(72, 14).Value = -22859.0763
$ws.Cells.Item(76, 8).Value = 4276481
$ws.Cells.Item(76, 10).Value = 6175762
$ws.Cells.Item(76, 12).Value = 6175762
$ws.Cells.Item(76, 14).Value = -6176392
$ws.Cells.Item(79, 8).Value = 4276481
$ws.Cells.Item(79, 10).Value = 6175762
$ws.Cells.Item(79, 12).Value = 6175762
$ws.Cells.Item(79, 14).Value = -6177946
$ws.Cells.Item(86, 8).Value = 5823.136
$ws.Cells.Item(86, 9).Value = 1204.5385
$ws.Cells.Item(86, 11).Value = 1204.5385
$ws.Cells.Item(86, 13).Value = -81.53850000000011
$ws.Cells.Item(89, 8).Value = 5823.136
$ws.Cells.Item(89, 9).Value = 1204.5385
$ws.Cells.Item(89, 11).Value = 6022.692500000001
$ws.Cells.Item(89, 13).Value = -406.692500000001
$ws.Cells.Item(129, 8).Value = 763.7593000000001
$ws.Cells.Item(129, 10).Value = 794.92
$ws.Cells.Item(129, 12).Value = 2384.76
$ws.Cells.Item(129, 14).Value = -12384.76
$ws.Cells.Item(132, 8).Value = 2578.7441
$ws.Cells.Item(132, 9).Value = 2625.9443
$ws.Cells.Item(132, 11).Value = 7877.8329
$ws.Cells.Item(132, 13).Value = -5347.8329
$ws.Cells.Item(138, 8).Value = 2249.1948
$ws.Cells.Item(138, 9).Value = 1172.0714
$ws.Cells.Item(138, 10).Value = 2864.6938
$ws.Cells.Item(138, 11).Value = 3516.2142
$ws.Cells.Item(138, 12).Value = 8594.081399999999
$ws.Cells.Item(138, 13).Value = 1623.7858
$ws.Cells.Item(138, 14).Value = -18874.0814

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6291.4062
$ws.Cells.Item(32, 9).Value = 4928.325
$ws.Cells.Item(32, 10).Value = 13106.8125
$ws.Cells.Item(32, 11).Value = 4928.325
$ws.Cells.Item(32, 12).Value = 13106.8125
$ws.Cells.Item(32, 13).Value = -4641.325
$ws.Cells.Item(32, 14).Value = -13680.8125
$ws.Cells.Item(63, 8).Value = 3908649.8
$ws.Cells.Item(63, 9).Value = 2742.5715
$ws.Cells.Item(63, 11).Value = 2742.5715
$ws.Cells.Item(63, 13).Value = -2056.5715
$ws.Cells.Item(66, 8).Value = 3908649.8
$ws.Cells.Item(66, 9).Value = 2742.5715
$ws.Cells.Item(66, 11).Value = 13712.8575
$ws.Cells.Item(66, 13).Value = -10280.8575
$ws.Cells.Item(97, 8).Value = 1224.75
$ws.Cells.Item(97, 9).Value = 1299.6666
$ws.Cells.Item(97, 11).Value = 1299.6666
$ws.Cells.Item(97, 13).Value = -803.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 989.6667
$ws.Cells.Item(94, 9).Value = 822.2727
$ws.Cells.Item(94, 10).Value = 1450
$ws.Cells.Item(94, 11).Value = 822.2727
$ws.Cells.Item(94, 12).Value = 1450
$ws.Cells.Item(94, 13).Value = -371.2727
$ws.Cells.Item(94, 14).Value = -2352
$ws.Cells.Item(96, 8).Value = 26574
$ws.Cells.Item(96, 9).Value = 10428
$ws.Cells.Item(96, 11).Value = 10428
$ws.Cells.Item(96, 13).Value = -7682
$ws.Cells.Item(105, 8).Value = 821510.9399999999
$ws.Cells.Item(105, 9).Value = 1376.1538
$ws.Cells.Item(105, 11).Value = 1376.1538
$ws.Cells.Item(105, 13).Value = 370.8462

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4325.7207
$ws.Cells.Item(31, 9).Value = 2253.389
$ws.Cells.Item(31, 10).Value = 5817.8
$ws.Cells.Item(31, 11).Value = 2253.389
$ws.Cells.Item(31, 12).Value = 5817.8
$ws.Cells.Item(31, 13).Value = -1958.389
$ws.Cells.Item(31, 14).Value = -6407.8
$ws.Cells.Item(34, 8).Value = 4325.7207
$ws.Cells.Item(34, 9).Value = 2253.389
$ws.Cells.Item(34, 10).Value = 5817.8
$ws.Cells.Item(34, 11).Value = 2253.389
$ws.Cells.Item(34, 12).Value = 5817.8
$ws.Cells.Item(34, 13).Value = -2051.389
$ws.Cells.Item(34, 14).Value = -6221.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 3120.7
$ws.Cells.Item(51, 9).Value = 1799.5
$ws.Cells.Item(51, 10).Value = 3451
$ws.Cells.Item(51, 11).Value = 5398.5
$ws.Cells.Item(51, 12).Value = 10353
$ws.Cells.Item(51, 13).Value = -4938.5
$ws.Cells.Item(51, 14).Value = -11273
$ws.Cells.Item(62, 8).Value = 4738.769
$ws.Cells.Item(62, 9).Value = 1847
$ws.Cells.Item(62, 10).Value = 9365.6
$ws.Cells.Item(62, 11).Value = 5541
$ws.Cells.Item(62, 12).Value = 28096.8
$ws.Cells.Item(62, 13).Value = -4855
$ws.Cells.Item(62, 14).Value = -29468.8
$ws.Cells.Item(63, 8).Value = 4241
$ws.Cells.Item(63, 9).Value = 2350
$ws.Cells.Item(63, 10).Value = 6132
$ws.Cells.Item(63, 11).Value = 7050
$ws.Cells.Item(63, 12).Value = 18396
$ws.Cells.Item(63, 13).Value = -6301
$ws.Cells.Item(63, 14).Value = -19894
$ws.Cells.Item(65, 8).Value = 4738.769
$ws.Cells.Item(65, 9).Value = 1847
$ws.Cells.Item(65, 10).Value = 9365.6
$ws.Cells.Item(65, 11).Value = 16623
$ws.Cells.Item(65, 12).Value = 84290.40000000001
$ws.Cells.Item(65, 13).Value = -13191
$ws.Cells.Item(65, 14).Value = -91154.40000000001
$ws.Cells.Item(66, 8).Value = 4241
$ws.Cells.Item(66, 9).Value = 2350
$ws.Cells.Item(66, 10).Value = 6132
$ws.Cells.Item(66, 11).Value = 21150
$ws.Cells.Item(66, 12).Value = 55188
$ws.Cells.Item(66, 13).Value = -17406
$ws.Cells.Item(66, 14).Value = -62676
$ws.Cells.Item(74, 8).Value = 9966.666999999999
$ws.Cells.Item(74, 10).Value = 9966.666999999999
$ws.Cells.Item(74, 12).Value = 29900.001
$ws.Cells.Item(74, 14).Value = -32022.001
$ws.Cells.Item(77, 8).Value = 9966.666999999999
$ws.Cells.Item(77, 10).Value = 9966.666999999999
$ws.Cells.Item(77, 12).Value = 89700.003
$ws.Cells.Item(77, 14).Value = -100308.003
$ws.Cells.Item(81, 8).Value = 5466.4
$ws.Cells.Item(81, 10).Value = 5466.4
$ws.Cells.Item(81, 12).Value = 16399.2
$ws.Cells.Item(81, 14).Value = -18645.2
$ws.Cells.Item(84, 8).Value = 5466.4
$ws.Cells.Item(84, 10).Value = 5466.4
$ws.Cells.Item(84, 12).Value = 49197.6
$ws.Cells.Item(84, 14).Value = -60429.6
$ws.Cells.Item(131, 8).Value = 884.9474
$ws.Cells.Item(131, 10).Value = 912.5571
$ws.Cells.Item(131, 12).Value = 2737.6713
$ws.Cells.Item(131, 14).Value = -12817.6713

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1425.7307
$ws.Cells.Item(97, 9).Value = 1498.7
$ws.Cells.Item(97, 10).Value = 1182.5
$ws.Cells.Item(97, 11).Value = 1498.7
$ws.Cells.Item(97, 12).Value = 1182.5
$ws.Cells.Item(97, 13).Value = -1002.7
$ws.Cells.Item(97, 14).Value = -2174.5
$ws.Cells.Item(122, 8).Value = 3508
$ws.Cells.Item(122, 9).Value = 2773.1667
$ws.Cells.Item(122, 11).Value = 8319.500100000001
$ws.Cells.Item(122, 13).Value = -5869.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2032.8572
$ws.Cells.Item(82, 9).Value = 2032.8572
$ws.Cells.Item(82, 11).Value = 2032.8572
$ws.Cells.Item(82, 13).Value = -1671.8572
$ws.Cells.Item(85, 8).Value = 2032.8572
$ws.Cells.Item(85, 9).Value = 2032.8572
$ws.Cells.Item(85, 11).Value = 2032.8572
$ws.Cells.Item(85, 13).Value = -784.8571999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2654.3333
$ws.Cells.Item(81, 9).Value = 633.3333
$ws.Cells.Item(81, 10).Value = 4675.3335
$ws.Cells.Item(81, 11).Value = 1266.6666
$ws.Cells.Item(81, 12).Value = 9350.666999999999
$ws.Cells.Item(81, 13).Value = -205.6666
$ws.Cells.Item(81, 14).Value = -11472.667
$ws.Cells.Item(84, 8).Value = 2654.3333
$ws.Cells.Item(84, 9).Value = 633.3333
$ws.Cells.Item(84, 10).Value = 4675.3335
$ws.Cells.Item(84, 11).Value = 6333.333000000001
$ws.Cells.Item(84, 12).Value = 46753.335
$ws.Cells.Item(84, 13).Value = -1029.333000000001
$ws.Cells.Item(84, 14).Value = -57361.335
$ws.Cells.Item(96, 8).Value = 2317.375
$ws.Cells.Item(96, 9).Value = 2243
$ws.Cells.Item(96, 10).Value = 2362
$ws.Cells.Item(96, 11).Value = 2243
$ws.Cells.Item(96, 12).Value = 2362
$ws.Cells.Item(96, 13).Value = -870
$ws.Cells.Item(96, 14).Value = -5108
$ws.Cells.Item(126, 8).Value = 2084.2
$ws.Cells.Item(126, 9).Value = 1695.4762
$ws.Cells.Item(126, 11).Value = 5086.4286
$ws.Cells.Item(126, 13).Value = -2616.4286
$ws.Cells.Item(133, 8).Value = 32215
$ws.Cells.Item(133, 10).Value = 32215
$ws.Cells.Item(133, 12).Value = 32215
$ws.Cells.Item(133, 14).Value = -42335
